$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text so numeric-looking
# strings (e.g. "0.998", "57.050.33") are preserved exactly as authored.
$ws.Range("D2:E51").NumberFormat = "@"

# Coin name / link swap + value updates
$ws.Range('D2').Value = '57.050.33'
$ws.Range('E2').Value = '  +3.95%  '
$ws.Range('D3').Value = '2.358.35'
$ws.Range('E3').Value = '  +2.75%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '521.22'
$ws.Range('E5').Value = '  +2.65%  '
$ws.Range('D6').Value = '135.06'
$ws.Range('E6').Value = '  +4.03%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '0.539'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').Value = '2.356.97'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('E10').Value = '  +6.75%  '
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('E12').Value = '  +4.74%  '
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').Value = '23.85'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = '2.777.48'
$ws.Range('E15').Value = '  +2.72%  '
$ws.Range('D16').Value = '56.995.54'
$ws.Range('E16').Value = '  +3.86%  '
$ws.Range('E17').Value = '  +2.52%  '
$ws.Range('D18').Value = '2.347.19'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').Value = '4.24'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').Value = '324.65'
$ws.Range('E21').Value = '  +4.84%  '
$ws.Range('D22').Value = '6.55'
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '61.08'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('E25').Value = '  +7.09%  '
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').Value = '7.88'
$ws.Range('E27').Value = '  +4.73%  '
$ws.Range('D28').Value = '1.28'
$ws.Range('E28').Value = '  +10.56%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = '170.93'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0749'
$ws.Range('E30').Value = '  +5.35%  '
$ws.Range('D31').Value = '1.70'
$ws.Range('E31').Value = '  +3.71%  '
$ws.Range('D32').Value = '6.20'
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').Value = '18.33'
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('E36').Value = '  +2.49%  '
$ws.Range('D37').Value = '0.925'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('D38').Value = '4.04'
$ws.Range('E38').Value = '  +4.16%  '
$ws.Range('D39').Value = '1.57'
$ws.Range('E39').Value = '  +8.98%  '
$ws.Range('E40').Value = '  +3.06%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('E42').Value = '  +4.35%  '
$ws.Range('D43').Value = '137.39'
$ws.Range('E43').Value = '  +2.18%  '
$ws.Range('D44').Value = '280.71'
$ws.Range('E44').Value = '  +9.58%  '
$ws.Range('E45').Value = '  +5.17%  '
$ws.Range('E46').Value = '  +2.60%  '
$ws.Range('D47').Value = '0.0506'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').Value = '0.565'
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('E49').Value = '  +4.67%  '
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '17.48'
$ws.Range('E51').Value = '  +5.46%  '
